$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D to Text format while writing the new values so that
# numeric-looking strings (e.g. "68.11", "7.20") keep their exact text
# representation instead of being auto-converted into numbers by Excel. The style
# is reset back to Normal afterwards so no extra formatting is left behind.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '58.376.50'
$ws.Range('E2').Value = '  -2.55%  '
$ws.Range('D3').Value = '3.139.72'
$ws.Range('E3').Value = '  -4.42%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '525.15'
$ws.Range('E5').Value = '  -5.00%  '
$ws.Range('D6').Value = '134.70'
$ws.Range('E6').Value = '  -3.98%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = '3.137.50'
$ws.Range('E8').Value = '  -4.58%  '
$ws.Range('E9').Value = '  -4.67%  '
$ws.Range('D10').Value = '7.20'
$ws.Range('E10').Value = '  -7.13%  '
$ws.Range('D11').Value = '0.108'
$ws.Range('E11').Value = '  -8.18%  '
$ws.Range('D12').Value = '0.380'
$ws.Range('E12').Value = '  -6.35%  '
$ws.Range('D13').Value = '3.678.82'
$ws.Range('E13').Value = '  -4.43%  '
$ws.Range('E14').Value = '  -1.08%  '
$ws.Range('D15').Value = '25.55'
$ws.Range('E15').Value = '  -4.91%  '
$ws.Range('D16').Value = '3.144.38'
$ws.Range('E16').Value = '  -4.17%  '
$ws.Range('D17').Value = '58.371.82'
$ws.Range('E17').Value = '  -2.74%  '
$ws.Range('E18').Value = '  -6.30%  '
$ws.Range('D19').Value = '5.76'
$ws.Range('E19').Value = '  -5.01%  '
$ws.Range('D20').Value = '13.01'
$ws.Range('E20').Value = '  -5.53%  '
$ws.Range('D21').Value = '7.91'
$ws.Range('E21').Value = '  -7.05%  '
$ws.Range('D22').Value = '344.68'
$ws.Range('E22').Value = '  -7.27%  '
$ws.Range('E23').Value = '  -0.12%  '
$ws.Range('D24').Value = '0.508'
$ws.Range('E24').Value = '  -4.32%  '
$ws.Range('D25').Value = '68.11'
$ws.Range('E25').Value = '  -7.32%  '
$ws.Range('D26').Value = '3.276.65'
$ws.Range('E26').Value = '  -4.19%  '
$ws.Range('D27').Value = '0.170'
$ws.Range('E27').Value = '  +0.35%  '
$ws.Range('D28').Value = '0.0₃0954'
$ws.Range('E28').Value = '  -5.71%  '
$ws.Range('E29').Value = '  -0.58%  '
$ws.Range('D30').Value = '6.81'
$ws.Range('E30').Value = '  -3.67%  '
$ws.Range('E31').Value = '  -0.08%  '
$ws.Range('D32').Value = '1.86'
$ws.Range('E32').Value = '  -7.78%  '
$ws.Range('D33').Value = '6.90'
$ws.Range('E33').Value = '  -7.39%  '
$ws.Range('E34').Value = '  -0.91%  '
$ws.Range('D35').Value = '21.41'
$ws.Range('E35').Value = '  -4.74%  '
$ws.Range('D36').Value = '4.81'
$ws.Range('E36').Value = '  -4.91%  '
$ws.Range('D37').Value = '157.39'
$ws.Range('E37').Value = '  -5.29%  '
$ws.Range('D38').Value = '6.23'
$ws.Range('E38').Value = '  -5.94%  '
$ws.Range('D39').Value = '1.37'
$ws.Range('E39').Value = '  -9.29%  '
$ws.Range('D40').Value = '0.0685'
$ws.Range('E40').Value = '  -5.41%  '
$ws.Range('D41').Value = '3.171.00'
$ws.Range('E41').Value = '  -4.43%  '
$ws.Range('B42').Value = 'OKB'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D42').Value = '40.41'
$ws.Range('E42').Value = '  -3.04%  '
$ws.Range('B43').Value = 'EnergySwap'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D43').Value = '24.23'
$ws.Range('E43').Value = '  -7.56%  '
$ws.Range('E44').Value = '  -1.20%  '
$ws.Range('D45').Value = '0.691'
$ws.Range('E45').Value = '  -7.23%  '
$ws.Range('D46').Value = '3.91'
$ws.Range('E46').Value = '  -4.72%  '
$ws.Range('E47').Value = '  -0.07%  '
$ws.Range('D48').Value = '1.44'
$ws.Range('E48').Value = '  -7.97%  '
$ws.Range('D49').Value = '2.274.95'
$ws.Range('E49').Value = '  -2.35%  '
$ws.Range('D50').Value = '6.20'
$ws.Range('E50').Value = '  -2.34%  '
$ws.Range('D51').Value = '20.77'
$ws.Range('E51').Value = '  -1.67%  '

$ws.Range("D2:D51").Style = "Normal"
